$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# misinterpreted as numbers (losing significant trailing zeros), so the
# stored text matches the source data exactly.
$textCells = @("D4","D6","D18","D21","D24","D36","D37","D47","D49")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '46.059.19'
$ws.Range("E2").Value = '  +1.09%  '
$ws.Range("D3").Value = '2.596.06'
$ws.Range("E3").Value = '  +7.38%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.21%  '
$ws.Range("D5").Value = '307.18'
$ws.Range("E5").Value = '  +4.45%  '
$ws.Range("D6").Value = '99.60'
$ws.Range("E6").Value = '  +5.40%  '
$ws.Range("D7").Value = '0.601'
$ws.Range("E7").Value = '  +5.94%  '
$ws.Range("E8").Value = '  +0.20%  '
$ws.Range("D9").Value = '0.579'
$ws.Range("E9").Value = '  +15.43%  '
$ws.Range("D10").Value = '39.22'
$ws.Range("E10").Value = '  +12.63%  '
$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").Value = '0.0844'
$ws.Range("E11").Value = '  +8.30%  '
$ws.Range("B12").Value = 'OKB'
$ws.Range("C12").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D12").Value = '54.15'
$ws.Range("E12").Value = '  +1.10%  '
$ws.Range("D13").Value = '8.18'
$ws.Range("E13").Value = '  +16.24%  '
$ws.Range("D14").Value = '3.004.95'
$ws.Range("E14").Value = '  +7.99%  '
$ws.Range("E15").Value = '  +1.50%  '
$ws.Range("D16").Value = '2.613.24'
$ws.Range("E16").Value = '  +8.01%  '
$ws.Range("D17").Value = '0.924'
$ws.Range("E17").Value = '  +10.26%  '
$ws.Range("D18").Value = '15.00'
$ws.Range("E18").Value = '  +6.26%  '
$ws.Range("D19").Value = '46.322.95'
$ws.Range("E19").Value = '  +2.09%  '
$ws.Range("E20").Value = '  +7.33%  '
$ws.Range("D21").Value = '13.00'
$ws.Range("E21").Value = '  +5.63%  '
$ws.Range("D22").Value = '6.72'
$ws.Range("E22").Value = '  +9.23%  '
$ws.Range("D23").Value = '71.37'
$ws.Range("E23").Value = '  +6.50%  '
$ws.Range("D24").Value = '272.70'
$ws.Range("E24").Value = '  +13.08%  '
$ws.Range("D25").Value = '3.03'
$ws.Range("E25").Value = '  +8.97%  '
$ws.Range("D26").Value = '29.96'
$ws.Range("E26").Value = '  +41.41%  '
$ws.Range("D27").Value = '2.17'
$ws.Range("E27").Value = '  +12.57%  '
$ws.Range("D28").Value = '0.999'
$ws.Range("E28").Value = '  -0.17%  '
$ws.Range("D29").Value = '4.01'
$ws.Range("E29").Value = '  +0.28%  '
$ws.Range("D30").Value = '10.58'
$ws.Range("E30").Value = '  +9.64%  '
$ws.Range("D31").Value = '2.32'
$ws.Range("E31").Value = '  +4.75%  '
$ws.Range("D32").Value = '38.93'
$ws.Range("E32").Value = '  +0.96%  '
$ws.Range("D33").Value = '6.25'
$ws.Range("E33").Value = '  +14.59%  '
$ws.Range("D34").Value = '3.64'
$ws.Range("E34").Value = '  -3.43%  '
$ws.Range("D35").Value = '2.83'
$ws.Range("E35").Value = '  +2.69%  '
$ws.Range("D36").Value = '0.0840'
$ws.Range("E36").Value = '  +9.57%  '
$ws.Range("D37").Value = '2.20'
$ws.Range("E37").Value = '  +11.04%  '
$ws.Range("D38").Value = '150.12'
$ws.Range("E38").Value = '  +1.27%  '
$ws.Range("D39").Value = '0.121'
$ws.Range("E39").Value = '  +6.84%  '
$ws.Range("E40").Value = '  +5.60%  '
$ws.Range("D41").Value = '23.08'
$ws.Range("E41").Value = '  +44.34%  '
$ws.Range("D42").Value = '15.88'
$ws.Range("E42").Value = '  +8.41%  '
$ws.Range("E44").Value = '  +10.59%  '
$ws.Range("D45").Value = '4.08'
$ws.Range("E45").Value = '  +8.31%  '
$ws.Range("D46").Value = '2.168.31'
$ws.Range("E46").Value = '  +8.76%  '
$ws.Range("D47").Value = '1.00'
$ws.Range("E47").Value = '  +0.16%  '
$ws.Range("D48").Value = '93.65'
$ws.Range("E48").Value = '  +5.59%  '
$ws.Range("D49").Value = '9.60'
$ws.Range("E49").Value = '  +12.78%  '
$ws.Range("D50").Value = '109.15'
$ws.Range("E50").Value = '  +8.59%  '
$ws.Range("E51").Value = '  -2.12%  '
